# SectorGroup.xlsx: re-order the "codeforiati" columns.
#
# Before: D=category-name, E=category-code, F=group-name,  G=group-code
# After:  D=group-code,    E=category-name, F=category-code, G=group-name
#
# i.e. the group-code column moves from G to D, and category-name/
# category-code/group-name each shift one column to the right
# (D->E, E->F, F->G) for every row, including the header row.
#
# We do the 4-column rotation via Range.Copy through a scratch area
# (columns Z:AC) so Excel preserves each cell's original type (text,
# stored as a shared string) instead of re-interpreting numeric-looking
# strings like "110" as numbers (which a plain .Value literal write
# would do, and which would also silently add a new cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 1
$lastRow = 235

$srcRange = "D" + $firstRow + ":G" + $lastRow
$scratchRange = "Z" + $firstRow

$ws.Range($srcRange).Copy($ws.Range($scratchRange))

$ws.Range("AC" + $firstRow + ":AC" + $lastRow).Copy($ws.Range("D" + $firstRow))
$ws.Range("Z" + $firstRow + ":Z" + $lastRow).Copy($ws.Range("E" + $firstRow))
$ws.Range("AA" + $firstRow + ":AA" + $lastRow).Copy($ws.Range("F" + $firstRow))
$ws.Range("AB" + $firstRow + ":AB" + $lastRow).Copy($ws.Range("G" + $firstRow))

$ws.Range("Z" + $firstRow + ":AC" + $lastRow).Clear()
